# Add new column 'Event' to Card24 by admin
# Net effect observed in the target XML:
#  - Column O ("servised by") is deleted entirely (header + all data cells),
#    shrinking the used range from A1:O12 to A1:N12.
#  - Column M ("Event") data cells (M2:M12) are cleared to empty values,
#    while the M1 header "Event" remains.
#  - Column N ("Correction") values are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Delete the entire column O ("servised by"). EntireColumn.Delete removes
# the column (header + all data) and updates the sheet's used
# range/dimension accordingly (A1:O12 -> A1:N12).
$ws.Range("O:O").EntireColumn.Delete()

# Clear the 'Event' column (M) data rows (M2:M12), leaving the M1 header
# "Event" intact.
$ws.Range("M2:M12").ClearContents()
